# Generate Report for Handoff
# The source UUID / commit-hash identifiers embedded in file names,
# hyperlink display text and the handoff timestamps are refreshed to a
# new report run.

$wb = $excel.ActiveWorkbook

$oldMd   = "60aad5fe-0fa1-48ea-9bbe-1c255fe73260.md"
$newMd   = "0061f551-f573-447e-9fa0-c4403e965861.md"

$oldZhXlf = "60aad5fe-0fa1-48ea-9bbe-1c255fe73260.c184118625916aaa0fc308aee916c53ea0a42ba7.zh-cn.xlf"
$newZhXlf = "0061f551-f573-447e-9fa0-c4403e965861.3bff11b70688e4a43285c2d0aa69fb00981ef098.zh-cn.xlf"
$newZhDateTime = "2016-03-04 11:03:05"

$oldDeXlf = "60aad5fe-0fa1-48ea-9bbe-1c255fe73260.c184118625916aaa0fc308aee916c53ea0a42ba7.de-de.xlf"
$newDeXlf = "0061f551-f573-447e-9fa0-c4403e965861.3bff11b70688e4a43285c2d0aa69fb00981ef098.de-de.xlf"
$newDeDateTime = "2016-03-04 11:03:20"

$configDisplay = ".localization-config"

# ---------------------------------------------------------------------------
# Overview sheet: only the .md hyperlink / display text changes.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/9d8254276edac3a7d6e53fa04ecde1b00c71895b/e2e/$oldMd"
$configAddress = "https://github.com/OpenLocalizationTest/oltest/blob/9d8254276edac3a7d6e53fa04ecde1b00c71895b/$configDisplay"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddress, [Type]::Missing, [Type]::Missing, $newMd)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configAddress, [Type]::Missing, [Type]::Missing, $configDisplay)

# ---------------------------------------------------------------------------
# zh-cn sheet: .md hyperlink, the zh-cn handoff xlf hyperlink + its
# "Latest Handoff Datetime" value.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d4cdf5aa42fb8d05e5881e977ed315d3eff00b82/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/ht/$oldZhXlf"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdAddress, [Type]::Missing, [Type]::Missing, $newMd)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), $zhXlfAddress, [Type]::Missing, [Type]::Missing, $newZhXlf)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $configAddress, [Type]::Missing, [Type]::Missing, $configDisplay)

$wsZhCn.Range("D2").Value2 = $newZhDateTime

# ---------------------------------------------------------------------------
# de-de sheet: .md hyperlink, the de-de handoff xlf hyperlink + its
# "Latest Handoff Datetime" value.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b111494b8a0d116793cd915aba3aaf0445b3d696/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/ht/$oldDeXlf"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdAddress, [Type]::Missing, [Type]::Missing, $newMd)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), $deXlfAddress, [Type]::Missing, [Type]::Missing, $newDeXlf)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $configAddress, [Type]::Missing, [Type]::Missing, $configDisplay)

$wsDeDe.Range("D2").Value2 = $newDeDateTime
